$d = $word.ActiveDocument

# The header table cell holds the name "Ben Bar" typed across three runs:
#   "Ben Ba" + "r" + "rrr"   ->  renders as "Ben Barrrr"
# Fix the typo by removing the trailing extra "rrr", leaving "Ben Bar".
# We locate the stray "rrr" via Find (so we are not dependent on hard-coded
# absolute character offsets) and delete just that trailing run's text with
# a Range.Delete(), which removes the whole now-empty run instead of merging
# the surrounding runs together.

$search = $d.Content
$found = $search.Find.Execute("Barrrr", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $delStart = $search.End - 3
    $delEnd = $search.End
    $extra = $d.Range($delStart, $delEnd)
    $extra.Delete()
}
